# Generate Report for Handback
# Populates the "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns on the zh-cn and de-de report sheets, and updates the "Status" column to
# reflect that the handback is now in sync with en-US.

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b5988bee95a256e4cc9e49f484317f1caa0ecae/e2e/"

$mdFile1 = "8c2dbc05-1860-4b0e-ace3-473864e757b5.md"
$mdFile2 = "9841f33d-2cd6-4222-adf0-a00e134293ba.md"

$statusText = "Handed back: in sync with en-US"

# Width conversion: raw OOXML column width = ColumnWidth + 5/6 (observed empirically),
# and ColumnWidth itself is quantized to steps of 1/6 by this runtime.
$wide30 = 30 - (5/6)
$wide40 = 40 - (5/6)

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- Status column (C) ---
    $ws.Range("C2").Value = $statusText
    $ws.Range("C3").Value = $statusText

    # --- Column widths ---
    $ws.Columns.Item(3).ColumnWidth = $wide30    # C: Status
    $ws.Columns.Item(9).ColumnWidth = $wide40    # I: Latest Target File
    $ws.Columns.Item(10).ColumnWidth = $wide40   # J: Latest Handback File

    # --- Row 2 (8c2dbc05...) ---
    $ws.Hyperlinks.Add($ws.Range("I2"), ($baseUrl + $mdFile1), "", "", $mdFile1)
    $ws.Range("I2").Style = "HyperLink"

    # --- Row 3 (9841f33d...) ---
    $ws.Hyperlinks.Add($ws.Range("I3"), ($baseUrl + $mdFile2), "", "", $mdFile2)
    $ws.Range("I3").Style = "HyperLink"
}

# --- zh-cn (sheet2) specific: Latest Handback File / DateTime ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("J2").Value = "8c2dbc05-1860-4b0e-ace3-473864e757b5.cc5abd98cf8c1a0eadd7798c362b19090aa5219c.zh-cn.xlf"
$wsZh.Range("J3").Value = "9841f33d-2cd6-4222-adf0-a00e134293ba.7c7cf04931d9eefa4d34884f02684cb0adc79975.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-03 16:35:31"
$wsZh.Range("K3").Value = "2016-09-03 16:35:31"

# --- de-de (sheet3) specific: Latest Handback File / DateTime ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("J2").Value = "8c2dbc05-1860-4b0e-ace3-473864e757b5.cc5abd98cf8c1a0eadd7798c362b19090aa5219c.de-de.xlf"
$wsDe.Range("J3").Value = "9841f33d-2cd6-4222-adf0-a00e134293ba.7c7cf04931d9eefa4d34884f02684cb0adc79975.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-03 16:35:38"
$wsDe.Range("K3").Value = "2016-09-03 16:35:38"

# --- Overview sheet column widths (E: zh-cn, F: de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $wide30
$wsOverview.Columns.Item(6).ColumnWidth = $wide30

Write-Host "Handback report generated."
